$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Size" table (rows 53-57): rename fields s_type/s_number -> s_unittype/s_value,
# give s_value the same VARCHAR(10) type as s_unittype, and drop the old
# s_description row (and its "M" value) entirely.
$ws.Rows("56:57").Delete()

$ws.Range("C55").Value = "s_unittype"
$ws.Range("C56").Value = "s_value"
$ws.Range("D56").Value = "VARCHAR(10)"

# Update the saved view/selection to match where the edit left off.
$ws.Activate()
[void]$ws.Range("F55").Select()
$excel.ActiveWindow.ScrollRow = 39
